$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 89, shifting existing rows 89:101 down to 90:102
$ws.Rows("89:89").Insert()

# Populate the new row 89 with this week's price report (same market/product
# as the surrounding rows, new date and price figures)
$ws.Range("A89").Value = 11
$ws.Range("B89").Value = "Vega Monumental Concepción"
$ws.Range("C89").Value = "Bíobío"
$ws.Range("D89").Value = "2021-10-05"
$ws.Range("E89").Value = 8
$ws.Range("F89").Value = "Fruta"
$ws.Range("G89").Value = 100108
$ws.Range("H89").Value = "Tropicales y subtropicales"
$ws.Range("I89").Value = 100108005
$ws.Range("J89").Value = "Piña"
$ws.Range("K89").Value = "Caramelo"
$ws.Range("L89").Value = "Segunda"
$ws.Range("M89").Value = 100
$ws.Range("N89").Value = 20000
$ws.Range("O89").Value = 21000
$ws.Range("P89").Value = 20500
$ws.Range("Q89").Value = "$/caja 14 unidades"
$ws.Range("R89").Value = "Ecuador"
$ws.Range("S89").Value = 1464
$ws.Range("T89").Value = 14
